$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -5
    9  = 2
    10 = 1
    11 = -10
    13 = 1
    14 = 4
    22 = -1
    25 = -4
    28 = -2
    31 = -9
    34 = -5
    40 = 0
    42 = 0
    43 = 5
    44 = -4
    50 = -4
    60 = 2
    61 = -6
    63 = 0
    64 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
